$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.087.11"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "2.320.55"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'541.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'131.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "2.319.41"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'0.332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'23.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "2.736.19"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "59.079.24"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "2.325.14"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "'10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "'313.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'6.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'62.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'0.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("D27").Value = "'7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").Value = "'171.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.31%  "
$ws.Range("D32").Value = "0.0₃0738"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").Value = "'1.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.22%  "
$ws.Range("D35").Value = "'0.384"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'17.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").Value = "'314.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.35%  "
$ws.Range("D41").Value = "'37.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'1.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").Value = "'3.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'0.0954"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").Value = "'18.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'0.0209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.26%  "

Write-Host "Applied cryptos update"
